$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row - Mikkel / raiseAlaram() (header-ish row right under the existing table,
# highlighted like the "Chris" rows below with a light fill)
$ws.Range("A9").Value = "Mikkel"
$ws.Range("B9").Value = "raiseAlaram()"

# New task row - Mikkel / deactivateAlarm()
$ws.Range("A10").Value = "Mikkel"
$ws.Range("B10").Value = "deactivateAlarm()"

# Light (background/white) fill highlight on the new "Mikkel" row
$ws.Range("A9:B9").Interior.ThemeColor = 2

# Leave the cursor on B10, matching where editing finished
$ws.Range("B10").Select() | Out-Null
